$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E: "Correction" header + ABS((C-D)/2) formula for rows 2-136 ---
$ws.Range("E1").Value = "Correction"
$ws.Range("F1").Value = "Title"
$ws.Range("E2").Formula = "=ABS((C2-D2)/2)"
$ws.Range("E3:E66").Formula = "=ABS((C3-D3)/2)"
$ws.Range("E67:E130").Formula = "=ABS((C67-D67)/2)"
$ws.Range("E131:E136").Formula = "=ABS((C131-D131)/2)"

# --- Remove stale STDEV formula from F5 ---
$ws.Range("F5").ClearContents()

# --- Column F: Title flag for rows 95-136 (2 = signal lost, 1 = detected) ---
$ws.Range("F95:F136").Value = 2
$ws.Range("F96").Value = 1

# --- Column E width ---
$ws.Columns.Item(5).ColumnWidth = 9.7

# --- AutoFilter over the data range ---
$ws.Range("A1:F136").AutoFilter()

# Excel auto-creates a hidden sheet-scoped _FilterDatabase name when AutoFilter is applied via the UI
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "='campi-10-2.5'!`$A`$1:`$F`$136", $false)
$fd.Visible = $false

# --- View state: L94 selected (matches author's saved view) ---
$ws.Range("L94").Select()
